$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "66.170.26"
Set-TextValue $ws.Range("E2") "  -1.59%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.778.08"
Set-TextValue $ws.Range("E3") "  +1.56%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.10%  "

# Row 5
Set-TextValue $ws.Range("D5") "407.86"
Set-TextValue $ws.Range("E5") "  -3.17%  "

# Row 6
Set-TextValue $ws.Range("D6") "131.90"
Set-TextValue $ws.Range("E6") "  +0.71%  "

# Row 7
Set-TextValue $ws.Range("D7") "3.766.61"
Set-TextValue $ws.Range("E7") "  +1.53%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.607"
Set-TextValue $ws.Range("E8") "  -5.87%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.999"
Set-TextValue $ws.Range("E9") "  -0.04%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.727"
Set-TextValue $ws.Range("E10") "  -5.13%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.167"
Set-TextValue $ws.Range("E11") "  -8.90%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0000362"
Set-TextValue $ws.Range("E12") "  -9.99%  "

# Row 13
Set-TextValue $ws.Range("D13") "40.80"
Set-TextValue $ws.Range("E13") "  -5.46%  "

# Row 14
Set-TextValue $ws.Range("D14") "4.350.19"
Set-TextValue $ws.Range("E14") "  +1.11%  "

# Row 15
Set-TextValue $ws.Range("D15") "9.83"
Set-TextValue $ws.Range("E15") "  -4.68%  "

# Row 16
Set-TextValue $ws.Range("D16") "14.72"
Set-TextValue $ws.Range("E16") "  +12.61%  "

# Row 17
Set-TextValue $ws.Range("E17") "  -1.48%  "

# Row 18
Set-TextValue $ws.Range("D18") "3.788.89"
Set-TextValue $ws.Range("E18") "  +2.16%  "

# Row 19
Set-TextValue $ws.Range("D19") "19.45"
Set-TextValue $ws.Range("E19") "  -6.18%  "

# Row 20
Set-TextValue $ws.Range("D20") "66.203.11"
Set-TextValue $ws.Range("E20") "  -1.57%  "

# Row 21
Set-TextValue $ws.Range("E21") "  -5.97%  "

# Row 22
Set-TextValue $ws.Range("D22") "411.35"
Set-TextValue $ws.Range("E22") "  -8.73%  "

# Row 23
Set-TextValue $ws.Range("D23") "14.36"
Set-TextValue $ws.Range("E23") "  -9.39%  "

# Row 24
Set-TextValue $ws.Range("D24") "85.24"
Set-TextValue $ws.Range("E24") "  -4.89%  "

# Row 25
Set-TextValue $ws.Range("D25") "3.05"
Set-TextValue $ws.Range("E25") "  -3.79%  "

# Row 26
Set-TextValue $ws.Range("D26") "5.71"
Set-TextValue $ws.Range("E26") "  +14.39%  "

# Row 27
Set-TextValue $ws.Range("D27") "35.91"
Set-TextValue $ws.Range("E27") "  -6.19%  "

# Row 28
Set-TextValue $ws.Range("D28") "3.12"
Set-TextValue $ws.Range("E28") "  -6.20%  "

# Row 29
Set-TextValue $ws.Range("E29") "  -8.72%  "

# Row 30
Set-TextValue $ws.Range("D30") "12.37"
Set-TextValue $ws.Range("E30") "  -1.75%  "

# Row 31
Set-TextValue $ws.Range("D31") "0.120"
Set-TextValue $ws.Range("E31") "  -1.58%  "

# Row 32
Set-TextValue $ws.Range("D32") "2.73"
Set-TextValue $ws.Range("E32") "  +0.75%  "

# Row 33
Set-TextValue $ws.Range("D33") "7.32"
Set-TextValue $ws.Range("E33") "  -0.63%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.154"
Set-TextValue $ws.Range("E34") "  -6.42%  "

# Row 35
Set-TextValue $ws.Range("D35") "38.94"
Set-TextValue $ws.Range("E35") "  -6.65%  "

# Row 36
Set-TextValue $ws.Range("E36") "  +0.13%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -3.01%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.0₃0753"
Set-TextValue $ws.Range("E38") "  -0.96%  "

# Row 39
Set-TextValue $ws.Range("E39") "  -7.22%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.85"
Set-TextValue $ws.Range("E40") "  -11.20%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.996"
Set-TextValue $ws.Range("E41") "  -0.25%  "

# Row 42
Set-TextValue $ws.Range("B42") "EnergySwap"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D42") "27.33"
Set-TextValue $ws.Range("E42") "  -3.10%  "

# Row 43
Set-TextValue $ws.Range("B43") "Stellar"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D43") "0.135"
Set-TextValue $ws.Range("E43") "  -8.81%  "

# Row 44
Set-TextValue $ws.Range("D44") "3.19"
Set-TextValue $ws.Range("E44") "  +20.10%  "

# Row 45
Set-TextValue $ws.Range("D45") "145.06"
Set-TextValue $ws.Range("E45") "  -1.99%  "

# Row 46
Set-TextValue $ws.Range("D46") "3.25"
Set-TextValue $ws.Range("E46") "  -5.20%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.06"
Set-TextValue $ws.Range("E47") "  -2.63%  "

# Row 48
Set-TextValue $ws.Range("E48") "  -2.66%  "

# Row 49
Set-TextValue $ws.Range("E49") "  -3.79%  "

# Row 50
Set-TextValue $ws.Range("D50") "2.79"
Set-TextValue $ws.Range("E50") "  -4.49%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.294"
Set-TextValue $ws.Range("E51") "  -4.90%  "
